# Add a new "2021" column (column R) to the water-samples table and
# update the current selection, matching the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (thin separator row above the header) ---
# R3 just needs the same (empty) formatting as its left neighbour Q3.
$ws.Range("Q3").Copy()
$ws.Range("R3").PasteSpecial(-4122)   # xlPasteFormats

# --- Row 4 (year header row) ---
# R4 gets the same formatting as Q4 plus the new year value.
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("R4").Value = 2021

# --- Row 5 (sanitary-chemical indicators) ---
# R5 keeps the same base formatting family as its row, but with a new
# vertical-top alignment variant (creates a fresh cellXf), plus the value.
$ws.Range("R5").Value = 0.9
$ws.Range("R5").VerticalAlignment = -4160   # xlTop

# --- Row 6 (microbiological indicators) ---
# R6 gets the same formatting as Q6 (0.0 number format) plus the value.
$ws.Range("Q6").Copy()
$ws.Range("R6").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("R6").Value = 6.5

# --- Update the active selection shown in the sheet view ---
$ws.Range("T5").Select()
